$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move active selection to E6 (per sheetView selection in the target)
$ws.Range("E6").Select()

# Update the F column "Sequences" values: replace placeholder "x" text with the
# actual number of sequences obtained per site, now that BLAST assignment is done.
$ws.Range("F6").Value = 4
$ws.Range("F8").Value = 4
$ws.Range("F12").Value = 4
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 4
$ws.Range("F16").Value = 3
$ws.Range("F18").Value = 4
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("F24").Value = 4
$ws.Range("F27").Value = 4
$ws.Range("F28").Value = 8
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 4
$ws.Range("F33").Value = 4

# Row 20 (Owen Park) site code corrected from "OP" to "Ow"
$ws.Range("E20").Value = "Ow"

# Remove the now-unused species BLAST-assignment helper codes (SM/SP/SR/SO)
# that used to live in column G.
$ws.Range("G20").Clear()
$ws.Range("G23").Clear()
$ws.Range("G25").Clear()
$ws.Range("G26").Clear()
